$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Theme 5 data import
#
# A new "WWMSP3" entry (WAMSI Westport Marine Science Program) needs
# review and is flagged in a new row (current row 35), highlighted in
# yellow, with a note in column K. The corrected entry "WWMSP3.1"
# (with a WQ Sensor site instead of WQ Grab) follows immediately
# after it, pushing the remaining UWA rows down by one.
# ------------------------------------------------------------------

# Insert a new row at 35 - this shifts the old rows 35-39 down to 36-40
# and copies the formatting (row height, cell styles) from row 34 above.
$ws.Rows("35:35").Insert()
$ws.Rows("35:35").RowHeight = 25.5

# --- Row 35 (new, highlighted "needs review" entry) -----------------
$ws.Range("A35").Value = "Data"
$ws.Range("B35").Value = "State Programs"
$ws.Range("C35").Value = "Western Australian Marine Science Institution"
$ws.Range("D35").Value = "WAMSI"
$ws.Range("E35").Value = "WAMSI Westport Marine Science Program"
$ws.Range("F35").Value = "WWMSP3"
$ws.Range("G35").Value = "WQ Grab"
$ws.Range("H35").Value = 18
$ws.Range("I35").Value = "Ongoing"
$ws.Range("J35").Value = "Y"

# --- Row 36 (corrected "WWMSP3.1" entry, replaces what used to be
#     row 35) ---------------------------------------------------------
$ws.Range("A36").ClearContents()
$ws.Range("B36").Value = "State Programs"
$ws.Range("C36").Value = "Western Australian Marine Science Institution"
$ws.Range("D36").Value = "WAMSI"
$ws.Range("E36").Value = "WAMSI Westport Marine Science Program"
$ws.Range("F36").Value = "WWMSP3.1"
$ws.Range("G36").Value = "WQ Sensor"
$ws.Range("H36").Value = 18
$ws.Range("I36").Value = "Ongoing"
$ws.Range("J36").ClearContents()

# The review note goes in last so the shared-string table gets
# "WWMSP3.1" before "Needs review, but should be removed".
$ws.Range("K35").Value = "Needs review, but should be removed"

# Highlight the new row (columns A-I) in yellow, preserving each
# column's existing alignment/wrap formatting.
$ws.Range("A35:I35").Interior.Color = 65535

# Reflect the reviewer's final view state (zoomed out a little, cursor
# left on the new note in K36).
$ws.Range("K36").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85

$wb.Save()
